# Sprint1 daily scrum meetings day 5 and 6
# Update task statuses and daily effort logging on the "Sprint" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")
# ($wb.ActiveSheet resolves to the same "Sprint" worksheet in this workbook;
# addressed by name here so the script is unambiguous either way.)

# "Create flow controls and buttons" task (row 9): moved from "In progress" to "Done",
# 1 hour of effort logged on Day 4 (column J).
$ws.Range("F9").Value = "Done"
$ws.Range("J9").Value = 1

# "Implement file chooser" task (row 10): moved from "In progress" to "Done",
# 2 hours of effort logged on Day 4 (column J).
$ws.Range("F10").Value = "Done"
$ws.Range("J10").Value = 2

# "Timeline" task (row 11): moved from "To do" to "In progress".
$ws.Range("F11").Value = "In progress"

# "Review code" task (row 12): moved from "To do" to "In progress".
$ws.Range("F12").Value = "In progress"

# Update the active selection to reflect where the user last clicked.
$ws.Range("F11").Select()

$wb.Save()
